$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B11:D12").NumberFormat = "@"
$ws.Range("B14:D14").NumberFormat = "@"

# Enterprises density (per 1000 people) - row 11
$ws.Range("B11").Value = "21.97"
$ws.Range("C11").Value = "10.68"
$ws.Range("D11").Value = "32.65"

# Employment (% of total) - row 12
$ws.Range("B12").Value = "11.64"
$ws.Range("C12").Value = "41.19"
$ws.Range("D12").Value = "52.83"

# Enterprises (% of total) - row 14
$ws.Range("B14").Value = "66.32"
$ws.Range("C14").Value = "32.23"
$ws.Range("D14").Value = "98.56"
